$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 "Modelo" - copy formatting (bold, border, alignment) from E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update slightly-changed numeric values in row 2
$ws.Range("B2").Value = 0.3638702225807679
$ws.Range("D2").Value = 0.468209296615762

# New data cell F2 with model description
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"

$excel.CutCopyMode = $false
